$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 214.44444
$ws.Range("I9").Value = 219.71428
$ws.Range("J9").Value = 196
$ws.Range("K9").Value = 219.71428
$ws.Range("L9").Value = 196
$ws.Range("M9").Value = -50.71428
$ws.Range("N9").Value = -534

$ws.Range("H11").Value = 21775512
$ws.Range("I11").Value = 21775512
$ws.Range("K11").Value = 21775512
$ws.Range("M11").Value = -21775372

$ws.Range("H40").Value = 5368.7144
$ws.Range("J40").Value = 2399.5
$ws.Range("L40").Value = 2399.5
$ws.Range("N40").Value = -2749.5

$ws.Range("H117").Value = 25024372
$ws.Range("J117").Value = 25024372
$ws.Range("L117").Value = 25024372
$ws.Range("N117").Value = -25033550

$ws.Range("H137").Value = 4944.569
$ws.Range("I137").Value = 2146.6155
$ws.Range("J137").Value = 14037.917
$ws.Range("K137").Value = 6439.8465
$ws.Range("L137").Value = 42113.751
$ws.Range("M137").Value = -3889.8465
$ws.Range("N137").Value = -47213.751

$ws.Range("H138").Value = 5147.846
$ws.Range("I138").Value = 6458.4443
$ws.Range("J138").Value = 2199
$ws.Range("K138").Value = 19375.3329
$ws.Range("L138").Value = 6597
$ws.Range("M138").Value = -14235.3329
$ws.Range("N138").Value = -16877

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H29").Value = 11165
$ws.Range("I29").Value = 6500
$ws.Range("J29").Value = 13497.5
$ws.Range("K29").Value = 6500
$ws.Range("L29").Value = 13497.5
$ws.Range("M29").Value = -6192
$ws.Range("N29").Value = -14113.5

$ws.Range("H32").Value = 175998.02
$ws.Range("I32").Value = 206037.58
$ws.Range("K32").Value = 206037.58
$ws.Range("M32").Value = -205750.58

$ws.Range("H36").Value = 2542.7144
$ws.Range("I36").Value = 2542.7144
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 2542.7144
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -2196.7144
$ws.Range("N36").ClearContents()

$ws.Range("H61").Value = 3174.6667
$ws.Range("I61").Value = 3028.68
$ws.Range("K61").Value = 3028.68
$ws.Range("M61").Value = -2816.68

$ws.Range("H97").Value = 680.5
$ws.Range("I97").Value = 598.23334
$ws.Range("K97").Value = 598.23334
$ws.Range("M97").Value = -102.23334

$ws.Range("H102").Value = 5354.3
$ws.Range("I102").Value = 4084.7856
$ws.Range("J102").Value = 8316.5
$ws.Range("K102").Value = 4084.7856
$ws.Range("L102").Value = 8316.5
$ws.Range("M102").Value = -2462.7856
$ws.Range("N102").Value = -11560.5

$ws.Range("H132").Value = 808359.6
$ws.Range("I132").Value = 863784.4399999999
$ws.Range("J132").Value = 4700
$ws.Range("K132").Value = 2591353.32
$ws.Range("L132").Value = 14100
$ws.Range("M132").Value = -2588823.32
$ws.Range("N132").Value = -19160

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 3174.6667
$ws.Range("I136").Value = 3028.68
$ws.Range("K136").Value = 9086.039999999999
$ws.Range("M136").Value = -6536.039999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 475.6
$ws.Range("I12").Value = 539.75
$ws.Range("K12").Value = 539.75
$ws.Range("M12").Value = -371.75

$ws.Range("H105").Value = 3023.1667
$ws.Range("I105").Value = 2260
$ws.Range("K105").Value = 2260
$ws.Range("M105").Value = -513

$ws.Range("H132").Value = 90514
$ws.Range("J132").Value = 90514
$ws.Range("L132").Value = 90514
$ws.Range("N132").Value = -100634

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3583.8286
$ws.Range("I31").Value = 3447.9
$ws.Range("K31").Value = 3447.9
$ws.Range("M31").Value = -3152.9

$ws.Range("H34").Value = 3583.8286
$ws.Range("I34").Value = 3447.9
$ws.Range("K34").Value = 3447.9
$ws.Range("M34").Value = -3245.9

$ws.Range("H134").Value = 1310.8889
$ws.Range("I134").Value = 974.75
$ws.Range("K134").Value = 2924.25
$ws.Range("M134").Value = -389.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1643.5
$ws.Range("I5").Value = 1148.8
$ws.Range("J5").Value = 2468
$ws.Range("K5").Value = 3446.4
$ws.Range("L5").Value = 7404
$ws.Range("M5").Value = -3334.4
$ws.Range("N5").Value = -7628

$ws.Range("H60").Value = 1502
$ws.Range("J60").Value = 2750
$ws.Range("L60").Value = 8250
$ws.Range("N60").Value = -8752

$ws.Range("H135").Value = 1643.5
$ws.Range("I135").Value = 1148.8
$ws.Range("J135").Value = 2468
$ws.Range("K135").Value = 10339.2
$ws.Range("L135").Value = 22212
$ws.Range("M135").Value = -7804.199999999999
$ws.Range("N135").Value = -27282

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 10007372
$ws.Range("J10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("N10").Value = -10338

$ws.Range("H11").Value = 6885117
$ws.Range("I11").Value = 8324642
$ws.Range("J11").Value = 167332.67
$ws.Range("K11").Value = 8324642
$ws.Range("L11").Value = 167332.67
$ws.Range("M11").Value = -8324503
$ws.Range("N11").Value = -167610.67

$ws.Range("H18").Value = 11833.333
$ws.Range("I18").Value = 9500
$ws.Range("J18").Value = 13000
$ws.Range("K18").Value = 9500
$ws.Range("L18").Value = 13000
$ws.Range("M18").Value = -9207
$ws.Range("N18").Value = -13586

$ws.Range("H70").Value = 34177.484
$ws.Range("I70").Value = 31354.818
$ws.Range("J70").Value = 38954.31
$ws.Range("K70").Value = 31354.818
$ws.Range("L70").Value = 38954.31
$ws.Range("M70").Value = -31084.818
$ws.Range("N70").Value = -39494.31

$ws.Range("H73").Value = 34177.484
$ws.Range("I73").Value = 31354.818
$ws.Range("J73").Value = 38954.31
$ws.Range("K73").Value = 31354.818
$ws.Range("L73").Value = 38954.31
$ws.Range("M73").Value = -30418.818
$ws.Range("N73").Value = -40826.31

$ws.Range("H80").Value = 2697.6667
$ws.Range("I80").Value = 2697.6667
$ws.Range("K80").Value = 2697.6667
$ws.Range("M80").Value = -1699.6667

$ws.Range("H83").Value = 2697.6667
$ws.Range("I83").Value = 2697.6667
$ws.Range("K83").Value = 13488.3335
$ws.Range("M83").Value = -8496.333500000001

$ws.Range("H132").Value = 5172.934
$ws.Range("I132").Value = 5183.619
$ws.Range("J132").Value = 5121.154
$ws.Range("K132").Value = 15550.857
$ws.Range("L132").Value = 15363.462
$ws.Range("M132").Value = -13020.857
$ws.Range("N132").Value = -20423.462

$ws.Range("H141").Value = 110214.5
$ws.Range("J141").Value = 110214.5
$ws.Range("L141").Value = 110214.5
$ws.Range("N141").Value = -120574.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 9500
$ws.Range("I23").Value = 9500
$ws.Range("K23").Value = 9500
$ws.Range("M23").Value = -9270

$ws.Range("H26").Value = 10400
$ws.Range("J26").Value = 10400
$ws.Range("L26").Value = 10400
$ws.Range("N26").Value = -10990

$ws.Range("H55").Value = 1188.9474
$ws.Range("I55").Value = 1031.9445
$ws.Range("J55").Value = 1330.25
$ws.Range("K55").Value = 1031.9445
$ws.Range("L55").Value = 1330.25
$ws.Range("M55").Value = -858.9445000000001
$ws.Range("N55").Value = -1676.25

$ws.Range("H82").Value = 2187.875
$ws.Range("J82").Value = 2907.2
$ws.Range("L82").Value = 2907.2
$ws.Range("N82").Value = -3629.2

$ws.Range("H85").Value = 2187.875
$ws.Range("J85").Value = 2907.2
$ws.Range("L85").Value = 2907.2
$ws.Range("N85").Value = -5403.2

$ws.Range("H100").Value = 4750
$ws.Range("I100").Value = 4000
$ws.Range("K100").Value = 4000
$ws.Range("M100").Value = -3459

$ws.Range("H132").Value = 2382.442
$ws.Range("I132").Value = 2157.1843
$ws.Range("J132").Value = 4094.4
$ws.Range("K132").Value = 6471.5529
$ws.Range("L132").Value = 12283.2
$ws.Range("M132").Value = -3941.5529
$ws.Range("N132").Value = -17343.2

$ws.Range("H133").Value = 89319.664
$ws.Range("J133").Value = 89319.664
$ws.Range("L133").Value = 89319.664
$ws.Range("N133").Value = -94379.664

$ws.Range("H136").Value = 1794.3158
$ws.Range("I136").Value = 1796.8334
$ws.Range("J136").Value = 1749
$ws.Range("K136").Value = 5390.5002
$ws.Range("L136").Value = 5247
$ws.Range("M136").Value = -2840.5002
$ws.Range("N136").Value = -10347

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 19591.572
$ws.Range("I52").Value = 7500
$ws.Range("K52").Value = 7500
$ws.Range("M52").Value = -7274

$ws.Range("H69").Value = 23999.5
$ws.Range("J69").Value = 23999.5
$ws.Range("L69").Value = 23999.5
$ws.Range("N69").Value = -25497.5

$ws.Range("H72").Value = 23999.5
$ws.Range("J72").Value = 23999.5
$ws.Range("L72").Value = 71998.5
$ws.Range("N72").Value = -79486.5

$ws.Range("H126").Value = 3949
$ws.Range("I126").Value = 3257
$ws.Range("K126").Value = 9771
$ws.Range("M126").Value = -7301

$ws.Range("H132").Value = 2767.8215
$ws.Range("I132").Value = 2523.95
$ws.Range("K132").Value = 7571.849999999999
$ws.Range("M132").Value = -5041.849999999999

$ws.Range("H136").Value = 1049.017
$ws.Range("I136").Value = 908.2292
$ws.Range("K136").Value = 2724.6876
$ws.Range("M136").Value = -174.6876000000002

$ws.Range("H141").Value = 99357.5
$ws.Range("J141").Value = 99357.5
$ws.Range("L141").Value = 99357.5
$ws.Range("N141").Value = -109717.5
